$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.607.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "'1.689.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'314.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'0.3897"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.14%  "
$ws.Range("D8").Value = "'0.4035"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'1.495"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").Value = "'1.003"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("D11").Value = "'52.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").Value = "'0.08756"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").Value = "'7.565"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.27%  "
$ws.Range("D14").Value = "'24.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.50%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.959"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.00001349"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("D17").Value = "'1.690.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "'98.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").Value = "'0.07108"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").Value = "'19.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("D21").Value = "'7.286"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.14%  "
$ws.Range("D22").Value = "'1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "'14.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "'24.596.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "'3.011"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.93%  "
$ws.Range("D26").Value = "'2.353"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "'22.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").Value = "'162.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "'8.826"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +15.72%  "
$ws.Range("D30").Value = "'136.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").Value = "'5.220"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("D32").Value = "'1.868.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").Value = "'0.08861"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("D34").Value = "'7.442"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.07%  "
$ws.Range("D35").Value = "'1.040"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.62%  "
$ws.Range("D36").Value = "'1.982"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.20%  "
$ws.Range("D37").Value = "'0.02918"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.16%  "
$ws.Range("D38").Value = "'0.2738"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").Value = "'10.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.11%  "
$ws.Range("D40").Value = "'0.09135"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").Value = "'14.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("D42").Value = "'0.7820"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.37%  "
$ws.Range("D43").Value = "'1.463"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "'16.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.72%  "
$ws.Range("D45").Value = "'0.7191"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("D46").Value = "'2.591"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "'4.194"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").Value = "'1.336"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("D50").Value = "'137.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").Value = "'90.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.83%  "
